$wb = $excel.ActiveWorkbook

# ALC!row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 303.23334
$ws.Range("J33").Value = 46.25
$ws.Range("L33").Value = 46.25
$ws.Range("N33").Value = -504.25

# ALC!row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 8718.484
$ws.Range("J112").Value = 8962.8125
$ws.Range("L112").Value = 26888.4375
$ws.Range("N112").Value = -29104.4375

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1375.695
$ws.Range("I137").Value = 1196.1316
$ws.Range("J137").Value = 1700.619
$ws.Range("K137").Value = 3588.3948
$ws.Range("L137").Value = 5101.857
$ws.Range("M137").Value = -1038.3948
$ws.Range("N137").Value = -10201.857

# ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1953.18
$ws.Range("I138").Value = 1359.4
$ws.Range("J138").Value = 2019.1555
$ws.Range("K138").Value = 4078.2
$ws.Range("L138").Value = 6057.4665
$ws.Range("M138").Value = 1061.8
$ws.Range("N138").Value = -16337.4665

# ALC!row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 49995
$ws.Range("J139").Value = 49995
$ws.Range("L139").Value = 49995
$ws.Range("N139").Value = -60275

# ALC!row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1820.0488
$ws.Range("I141").Value = 1013.2564
$ws.Range("J141").Value = 17552.5
$ws.Range("K141").Value = 3039.7692
$ws.Range("L141").Value = 52657.5
$ws.Range("M141").Value = 2140.2308
$ws.Range("N141").Value = -63017.5

# ARM!row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 174
$ws.Range("I4").Value = 198
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 198
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -82
$ws.Range("N4").Value = -382

# ARM!row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 212.66667
$ws.Range("I5").Value = 180
$ws.Range("J5").Value = 219.2
$ws.Range("K5").Value = 180
$ws.Range("L5").Value = 219.2
$ws.Range("M5").Value = -68
$ws.Range("N5").Value = -443.2

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 680323.1
$ws.Range("I32").Value = 754430.25
$ws.Range("J32").Value = 20769.7
$ws.Range("K32").Value = 754430.25
$ws.Range("L32").Value = 20769.7
$ws.Range("M32").Value = -754143.25
$ws.Range("N32").Value = -21343.7

# ARM!row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2058.54
$ws.Range("I74").Value = 981.43475
$ws.Range("J74").Value = 2976.074
$ws.Range("K74").Value = 981.43475
$ws.Range("L74").Value = 2976.074
$ws.Range("M74").Value = -107.43475
$ws.Range("N74").Value = -4724.074000000001

# ARM!row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2058.54
$ws.Range("I77").Value = 981.43475
$ws.Range("J77").Value = 2976.074
$ws.Range("K77").Value = 4907.17375
$ws.Range("L77").Value = 14880.37
$ws.Range("M77").Value = -539.1737499999999
$ws.Range("N77").Value = -23616.37

# BSM!row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 212.66667
$ws.Range("I4").Value = 180
$ws.Range("J4").Value = 219.2
$ws.Range("K4").Value = 180
$ws.Range("L4").Value = 219.2
$ws.Range("M4").Value = -65
$ws.Range("N4").Value = -449.2

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1452
$ws.Range("I107").Value = 1145.7142
$ws.Range("J107").Value = 2166.6667
$ws.Range("K107").Value = 1145.7142
$ws.Range("L107").Value = 2166.6667
$ws.Range("M107").Value = 774.2858000000001
$ws.Range("N107").Value = -6006.6667

# CRP!row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 555623.8
$ws.Range("I7").Value = 694479.75
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 694479.75
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -694366.75
$ws.Range("N7").Value = -426

# CRP!row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1823
$ws.Range("I16").Value = 1690
$ws.Range("J16").Value = 1956
$ws.Range("K16").Value = 1690
$ws.Range("L16").Value = 1956
$ws.Range("M16").Value = -1403
$ws.Range("N16").Value = -2530

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4022.1143
$ws.Range("I31").Value = 1275.9744
$ws.Range("J31").Value = 7476.9355
$ws.Range("K31").Value = 1275.9744
$ws.Range("L31").Value = 7476.9355
$ws.Range("M31").Value = -980.9744000000001
$ws.Range("N31").Value = -8066.9355

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4022.1143
$ws.Range("I34").Value = 1275.9744
$ws.Range("J34").Value = 7476.9355
$ws.Range("K34").Value = 1275.9744
$ws.Range("L34").Value = 7476.9355
$ws.Range("M34").Value = -1073.9744
$ws.Range("N34").Value = -7880.9355

# CRP!row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1067.4166
$ws.Range("I105").Value = 1076.125
$ws.Range("J105").Value = 1050
$ws.Range("K105").Value = 1076.125
$ws.Range("L105").Value = 1050
$ws.Range("M105").Value = 670.875
$ws.Range("N105").Value = -4544

# CRP!row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2977046.8
$ws.Range("I107").Value = 12500598
$ws.Range("J107").Value = 937
$ws.Range("K107").Value = 12500598
$ws.Range("L107").Value = 937
$ws.Range("M107").Value = -12498678
$ws.Range("N107").Value = -4777

# CRP!row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1823
$ws.Range("I113").Value = 1690
$ws.Range("J113").Value = 1956
$ws.Range("K113").Value = 1690
$ws.Range("L113").Value = 1956
$ws.Range("M113").Value = 480
$ws.Range("N113").Value = -6296

# CUL!row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 9082.739
$ws.Range("J26").Value = 11583.333
$ws.Range("L26").Value = 34749.999
$ws.Range("N26").Value = -35325.999

# CUL!row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2455.8333
$ws.Range("I122").Value = 383.19354
$ws.Range("K122").Value = 3448.74186
$ws.Range("M122").Value = -998.7418600000001

# LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 93577.27
$ws.Range("I40").Value = 127425
$ws.Range("J40").Value = 3316.6667
$ws.Range("K40").Value = 127425
$ws.Range("L40").Value = 3316.6667
$ws.Range("M40").Value = -127289
$ws.Range("N40").Value = -3588.6667

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 705.5833
$ws.Range("I46").Value = 571.6667
$ws.Range("J46").Value = 1107.3334
$ws.Range("K46").Value = 571.6667
$ws.Range("L46").Value = 1107.3334
$ws.Range("M46").Value = -383.6667
$ws.Range("N46").Value = -1483.3334

# WVR!row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 960.5625
$ws.Range("I107").Value = 905.5
$ws.Range("J107").Value = 1125.75
$ws.Range("K107").Value = 2716.5
$ws.Range("L107").Value = 3377.25
$ws.Range("M107").Value = -796.5
$ws.Range("N107").Value = -7217.25

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9117216
$ws.Range("I132").Value = 3327.8235
$ws.Range("J132").Value = 19446290
$ws.Range("K132").Value = 9983.470499999999
$ws.Range("L132").Value = 58338870
$ws.Range("M132").Value = -7453.470499999999
$ws.Range("N132").Value = -58343930

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1838.68
$ws.Range("I136").Value = 1909.9342
$ws.Range("J136").Value = 1613.0416
$ws.Range("K136").Value = 5729.8026
$ws.Range("L136").Value = 4839.1248
$ws.Range("M136").Value = -3179.8026
$ws.Range("N136").Value = -9939.1248
